$p = $ppt.ActivePresentation

foreach ($idx in 14,15,16) {
    $s = $p.Slides.Item($idx)
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{F7AED1DB-CF54-4B7F-8904-104769ADFB61}")
        }
    }
}
